# Checkpoint 3 of the "Analysing the Top 3 English-Speaking Countries"
# table (sheet "Table-3.1"): fill in the answers for the top three
# English-speaking countries, then leave that sheet active/selected
# (mirrors the author having just finished working on it).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table-3.1")

$ws.Range("C5").Value = "USA"
$ws.Range("C6").Value = "GBR"
$ws.Range("C7").Value = "CAN"

$ws.Activate() | Out-Null
$ws.Range("B8").Select() | Out-Null
